$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2, 4, 6, 9, 11, 13, 15, 16, 18, 19, 20, 21, 22, 23)
foreach ($r in $rows) {
    $ws.Range("C$r").Value = "nan"
}
